$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.005899579586891775
$ws.Range("J2").Value = 0.008823342375055644
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.0108025
$ws.Range("N2").Value = 0.021605
$ws.Range("Q2").Value = 0.0006944567166666667
$ws.Range("R2").Value = 0.0041667403
$ws.Range("S2").Value = 0.005899579586891775
$ws.Range("T2").Value = 0.008823342375055644

# Row 3 updates
$ws.Range("G3").Value = 10.832535
$ws.Range("H3").Value = 21.66507
$ws.Range("I3").Value = 0.9941004204131083
$ws.Range("J3").Value = 0.9911766576249443
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.0108025
$ws.Range("N3").Value = 0.021605
$ws.Range("Q3").Value = 0.1170184593375
$ws.Range("R3").Value = 0.46807383735
$ws.Range("S3").Value = 0.9941004204131083
$ws.Range("T3").Value = 0.9911766576249443
